$wb = $excel.ActiveWorkbook

# --- 1. hotel_info: insert a new "State" column between Hotel_Name (B) and City (C) ---
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$hotelInfo.Range("C1").EntireColumn.Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"

# --- 2. Reorder worksheets so review_info comes before hotel_info ---
$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($hotelInfo)
